$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = '30.413.94'
$ws.Range("E2").Value = '  -0.04%  '

# Row 3
$ws.Range("D3").Value = '1.927.96'
$ws.Range("E3").Value = '  +4.16%  '

# Row 4
$ws.Range("D4").Value = '''0.9990'
$ws.Range("E4").Value = '  -0.12%  '

# Row 5
$ws.Range("D5").Value = '''240.64'
$ws.Range("E5").Value = '  +3.20%  '

# Row 6
$ws.Range("D6").Value = '''0.9995'
$ws.Range("E6").Value = '  -0.08%  '

# Row 7
$ws.Range("D7").Value = '''0.4752'
$ws.Range("E7").Value = '  -0.12%  '

# Row 8
$ws.Range("D8").Value = '''44.35'
$ws.Range("E8").Value = '  +2.54%  '

# Row 9
$ws.Range("D9").Value = '''0.2867'
$ws.Range("E9").Value = '  +4.37%  '

# Row 10
$ws.Range("D10").Value = '''0.06565'
$ws.Range("E10").Value = '  +3.83%  '

# Row 11
$ws.Range("D11").Value = '''19.04'
$ws.Range("E11").Value = '  +8.41%  '

# Row 12
$ws.Range("D12").Value = '''107.02'
$ws.Range("E12").Value = '  +26.33%  '

# Row 13
$ws.Range("D13").Value = '1.926.40'
$ws.Range("E13").Value = '  +3.99%  '

# Row 14
$ws.Range("D14").Value = '''0.07617'
$ws.Range("E14").Value = '  +2.12%  '

# Row 15
$ws.Range("D15").Value = '''5.130'
$ws.Range("E15").Value = '  +3.68%  '

# Row 16
$ws.Range("D16").Value = '''0.6561'
$ws.Range("E16").Value = '  +5.09%  '

# Row 17
$ws.Range("D17").Value = '''305.51'
$ws.Range("E17").Value = '  +23.91%  '

# Row 18
$ws.Range("D18").Value = '30.422.07'

# Row 19
$ws.Range("D19").Value = '''1.000'

# Row 20
$ws.Range("D20").Value = '''12.95'
$ws.Range("E20").Value = '  +2.40%  '

# Row 21
$ws.Range("D21").Value = '2.171.71'
$ws.Range("E21").Value = '  +3.86%  '

# Row 22
$ws.Range("D22").Value = '''0.000007488'
$ws.Range("E22").Value = '  +2.31%  '

# Row 23
$ws.Range("D23").Value = '''5.299'
$ws.Range("E23").Value = '  +8.09%  '

# Row 24
$ws.Range("D24").Value = '''0.9998'
$ws.Range("E24").Value = '  -0.08%  '

# Row 25
$ws.Range("D25").Value = '''6.257'
$ws.Range("E25").Value = '  +5.98%  '

# Row 26
$ws.Range("D26").Value = '''167.30'
$ws.Range("E26").Value = '  +1.65%  '

# Row 27
$ws.Range("D27").Value = '''9.214'
$ws.Range("E27").Value = '  +1.41%  '

# Row 28
$ws.Range("D28").Value = '''20.16'
$ws.Range("E28").Value = '  +12.16%  '

# Row 29
$ws.Range("D29").Value = '''2.025'
$ws.Range("E29").Value = '  +8.29%  '

# Row 30
$ws.Range("D30").Value = '''0.1110'
$ws.Range("E30").Value = '  +7.90%  '

# Row 31
$ws.Range("D31").Value = '''1.356'
$ws.Range("E31").Value = '  +0.65%  '

# Row 32
$ws.Range("D32").Value = '''4.082'
$ws.Range("E32").Value = '  +1.26%  '

# Row 33
$ws.Range("D33").Value = '''3.918'
$ws.Range("E33").Value = '  +2.62%  '

# Row 34
$ws.Range("D34").Value = '''0.04992'
$ws.Range("E34").Value = '  +3.16%  '

# Row 35
$ws.Range("D35").Value = '''0.7410'
$ws.Range("E35").Value = '  +6.44%  '

# Row 36
$ws.Range("D36").Value = '''1.147'
$ws.Range("E36").Value = '  +1.61%  '

# Row 37
$ws.Range("D37").Value = '''2.750'
$ws.Range("E37").Value = '  +1.77%  '

# Row 38
$ws.Range("B38").Value = 'Frax'
$ws.Range("C38").Value = 'https://coinranking.com/coin/KfWtaeV1W+frax-frax'
$ws.Range("D38").Value = '''0.9988'
$ws.Range("E38").Value = '  -0.08%  '

# Row 39
$ws.Range("B39").Value = 'VeChain'
$ws.Range("C39").Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range("D39").Value = '''0.01938'
$ws.Range("E39").Value = '  +2.01%  '

# Row 40
$ws.Range("B40").Value = 'MXToken'
$ws.Range("C40").Value = 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
$ws.Range("D40").Value = '''2.696'
$ws.Range("E40").Value = '  +0.63%  '

# Row 41
$ws.Range("B41").Value = 'RenderToken'
$ws.Range("C41").Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range("D41").Value = '''2.053'
$ws.Range("E41").Value = '  +2.90%  '

# Row 42
$ws.Range("B42").Value = 'TrustWalletToken'
$ws.Range("C42").Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$ws.Range("D42").Value = '''0.8770'
$ws.Range("E42").Value = '  -0.04%  '

# Row 43
$ws.Range("B43").Value = 'Quant'
$ws.Range("C43").Value = 'https://coinranking.com/coin/bauj_21eYVwso+quant-qnt'
$ws.Range("D43").Value = '''106.89'
$ws.Range("E43").Value = '  +0.15%  '

# Row 44
$ws.Range("B44").Value = 'FraxShare'
$ws.Range("C44").Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range("D44").Value = '''5.797'
$ws.Range("E44").Value = '  +5.42%  '

# Row 45
$ws.Range("B45").Value = 'Aave'
$ws.Range("C45").Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range("D45").Value = '''70.15'
$ws.Range("E45").Value = '  +11.16%  '

# Row 46
$ws.Range("B46").Value = 'PaxDollar'
$ws.Range("C46").Value = 'https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp'
$ws.Range("D46").Value = '''0.9996'
$ws.Range("E46").Value = '  -0.06%  '

# Row 47
$ws.Range("B47").Value = 'TheSandbox'
$ws.Range("C47").Value = 'https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand'
$ws.Range("D47").Value = '''0.4139'
$ws.Range("E47").Value = '  +2.17%  '

# Row 48
$ws.Range("B48").Value = 'Aptos'
$ws.Range("C48").Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$ws.Range("D48").Value = '''7.255'
$ws.Range("E48").Value = '  +1.15%  '

# Row 49
$ws.Range("B49").Value = 'EnergySwap'
$ws.Range("C49").Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range("D49").Value = '''9.271'
$ws.Range("E49").Value = '  +8.78%  '

# Row 50
$ws.Range("B50").Value = 'Elrond'
$ws.Range("C50").Value = 'https://coinranking.com/coin/omwkOTglq+elrond-egld'
$ws.Range("D50").Value = '''34.79'
$ws.Range("E50").Value = '  +3.30%  '

# Row 51
$ws.Range("B51").Value = 'Algorand'
$ws.Range("C51").Value = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
$ws.Range("D51").Value = '''0.1201'
$ws.Range("E51").Value = '  +0.34%  '
